# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-13
$newK = @{
    2  = 3
    3  = 7
    4  = 6
    5  = 4
    6  = 6
    7  = 12
    8  = 9
    9  = 3
    10 = 8
    11 = 5
    12 = 4
    13 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
